# Functional traits all recorded species.xlsx - workbook edit
#
# Summary of the change (per commit message / diff):
#  - Two species records were removed from Sheet1:
#       SpecID 934  (Dendrocopos moluccensis / Sunda Woodpecker)
#       SpecID 2878 (Tyto alba / Barn Owl)
#  - The remaining data rows were re-sorted by column I (English name)
#    instead of column A (SpecID) - this is reflected in the autofilter's
#    recorded sort state in the underlying XML.
#  - The active selection in the bottom-right (frozen) pane moved to F26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two retired species rows -------------------------------
# Delete from the bottom up so row indices above the deleted row aren't
# invalidated by the earlier delete.
$ws.Rows(9).Delete()   # SpecID 2878 - Barn Owl         (row 9 before any delete)
$ws.Rows(3).Delete()   # SpecID 934  - Sunda Woodpecker (row 3 before any delete)

# --- Re-sort the remaining 26 data rows (rows 2:27) by English name ----
$sortRange = $ws.Range("A1:AN27")
$key1 = $ws.Range("I1:I27")

$af = $ws.AutoFilter
$s = $af.Sort
$sf = $s.SortFields
$sf.Clear()
$sf.Add($key1)
$s.SetRange($sortRange)
$s.Header = 1
$s.Apply()

# --- Update the saved selection to match the author's final cursor pos -
$ws.Range("F26").Select()
